$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price / volume data.
# For cells whose new text would otherwise be auto-parsed as a number by Excel
# (single-dot decimals like "131.52"), briefly force a text number format so the
# value is written as a string, then restore the original (default) cell style.

$ws.Range("D2").Value = "42.130.81"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.265.32"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "131.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +13,040.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "93.48"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0804"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "2.621.33"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "2.269.30"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.786"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("D19").Value = "41.980.56"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.03%  "
$ws.Range("D21").Value = "0.0₃0918"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "159.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.38%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0745"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.97%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.009.86"
$ws.Range("E43").Value = "  -2.83%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("E45").Value = "  +10.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0283"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "
